# This script applies a cyclic rotation of columns D (group-code /
# codeforiati:group-code), E (category-name / codeforiati:category-name)
# and G (group-name / codeforiati:group-name) for every row of the sheet
# (including the header row). Column F (category-code) is left untouched.
#
# The rotation is: new D = old E, new E = old G, new G = old D.
#
# This matches the observed diff: header D/E/G text rotates from
#   D=group-code, E=category-name, G=group-name
# to
#   D=category-name, E=group-name, G=group-code
# and every data row's displayed D/E/G values rotate the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $gVal = $ws.Cells.Item($r, 7).Value2

    $ws.Cells.Item($r, 4).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $gVal
    $ws.Cells.Item($r, 7).Value = $dVal
}
